$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: PM_IPA_FERMENTACION_PRESION sample, valor 2.1 -> 7.45 ---
$ws.Range("B2").Value = 7.45

# --- Row 3: PM_IPA_FERMENTACION_PRESION sample, valor 2.1 -> 7.45 ---
$ws.Range("B3").Value = 7.45

# --- Row 4: PM_IPA_FERMENTACION_PRESION sample, valor 2.1 -> 7.45 ---
$ws.Range("B4").Value = 7.45

# --- Row 5: was PM_IPA_CENTRIFUGADO_MARCHA, becomes another PM_IPA_FERMENTACION_PRESION row ---
$ws.Range("A5").Value = "PM_IPA_FERMENTACION_PRESION"
$ws.Range("B5").Value = 7.45
$ws.Range("C5").Value = 6.8
$ws.Range("D5").Formula = "=IF(AND(B5>3,B5<7),""presion alta"",""presion normal"")"
$ws.Range("F5").Value = "text"

# --- Row 6: formula changes to simple green/blue comparison ---
$ws.Range("D6").Formula = "=IF(B6=C6,""green"",""blue"")"

# --- Row 7: formula changes to marcha combinada / marcha no combinada ---
$ws.Range("D7").Formula = "=IF(AND(B7>3,B7<7),""marcha combinada"",""marcha no combinada"")"

# --- Update active selection to F2 ---
$ws.Range("F2").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
"done"
